# Added Search customer by Name test case
# Update the Customer_info sheet's email column with corrected addresses
# (matching each customer's actual first+last name), editing row by row
# K3 -> K4 -> K5 -> K6 -> K2, finishing with K2 selected (as last edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer_info")

$ws.Range("K3").Value = "manikapoor@gmail.com"
$ws.Range("K4").Value = "meenasharma1@gmail.com"
$ws.Range("K5").Value = "saumyachopra@gmail.com"
$ws.Range("K6").Value = "vihangupta1@gmail.com"
$ws.Range("K2").Value = "rameshgoyal@gmail.com"

$null = $ws.Range("K2").Select()
